$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the measured values (cells C2, C3, B4)
$ws.Range("C2").Value = 10.5
$ws.Range("C3").Value = 9.5
$ws.Range("B4").Value = 0.85
